$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5538
$ws.Range("I3").Value = 7492
$ws.Range("J3").Value = 5916
$ws.Range("H4").Value = 1702
$ws.Range("I4").Value = 1773
$ws.Range("J4").Value = 1285
$ws.Range("J5").Value = 457
$ws.Range("J6").Value = 7454
$ws.Range("H7").Value = 26013
$ws.Range("I7").Value = 26228
$ws.Range("J7").Value = 20650

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 87
$ws.Range("J6").Value = 151
$ws.Range("J7").Value = 605
$ws.Range("J8").Value = 1297
$ws.Range("J9").Value = 99
$ws.Range("J10").Value = 139
$ws.Range("J11").Value = 321
$ws.Range("J13").Value = 26
$ws.Range("J18").Value = 175
$ws.Range("J19").Value = 603
$ws.Range("J20").Value = 424
$ws.Range("J22").Value = 54
$ws.Range("J23").Value = 197
$ws.Range("J26").Value = 45
$ws.Range("J27").Value = 124
$ws.Range("J29").Value = 1156
$ws.Range("J33").Value = 947
$ws.Range("J35").Value = 30
$ws.Range("J36").Value = 287
$ws.Range("J37").Value = 640
$ws.Range("J42").Value = 855
$ws.Range("J43").Value = 170
$ws.Range("J48").Value = 240
$ws.Range("J49").Value = 139
$ws.Range("J51").Value = 255
$ws.Range("J52").Value = 523
$ws.Range("J53").Value = 283
$ws.Range("J54").Value = 402
$ws.Range("I63").Value = 240
$ws.Range("J63").Value = 76
$ws.Range("J65").Value = 529
$ws.Range("J66").Value = 66
$ws.Range("J67").Value = 781
$ws.Range("J72").Value = 84
$ws.Range("J74").Value = 22
$ws.Range("J78").Value = 254
$ws.Range("J79").Value = 590
$ws.Range("J83").Value = 423
$ws.Range("J85").Value = 871
$ws.Range("J89").Value = 271
$ws.Range("I90").Value = 338
$ws.Range("J90").Value = 225
$ws.Range("J91").Value = 228
$ws.Range("H97").Value = 211
$ws.Range("J98").Value = 145
$ws.Range("J99").Value = 323
$ws.Range("J100").Value = 40
$ws.Range("H101").Value = 26013
$ws.Range("I101").Value = 26228
$ws.Range("J101").Value = 20650

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 187
$ws.Range("J3").Value = 182
$ws.Range("J6").Value = 196
$ws.Range("J7").Value = 605

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 98
$ws.Range("J3").Value = 66
$ws.Range("J6").Value = 131
$ws.Range("J7").Value = 321

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 85
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 271

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J4").Value = 57
$ws.Range("J5").Value = 17
$ws.Range("J7").Value = 871

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 161
$ws.Range("J7").Value = 523

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 180
$ws.Range("J7").Value = 283

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 357
$ws.Range("J3").Value = 395
$ws.Range("J6").Value = 436
$ws.Range("J7").Value = 1297

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 126
$ws.Range("J3").Value = 157
$ws.Range("J6").Value = 115
$ws.Range("J7").Value = 423

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 234
$ws.Range("J3").Value = 308
$ws.Range("J7").Value = 947

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 218
$ws.Range("J6").Value = 188
$ws.Range("J7").Value = 640

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 156
$ws.Range("J3").Value = 151
$ws.Range("J6").Value = 186
$ws.Range("J7").Value = 529

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J6").Value = 82
$ws.Range("J7").Value = 323

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 301
$ws.Range("J7").Value = 781

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J2").Value = 24
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 139

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 79
$ws.Range("J6").Value = 190
$ws.Range("J7").Value = 402

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 345
$ws.Range("J3").Value = 399
$ws.Range("J7").Value = 1156

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J4").Value = 37
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 240

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 177
$ws.Range("J6").Value = 222
$ws.Range("J7").Value = 603

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 151

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J4").Value = 40
$ws.Range("J6").Value = 438
$ws.Range("J7").Value = 855

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("J5").Value = 12
$ws.Range("J6").Value = 26

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 139

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J3").Value = 84
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 254

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J2").Value = 52
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 197

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 228

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J5").Value = 17
$ws.Range("J6").Value = 165
$ws.Range("J7").Value = 590

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J3").Value = 148
$ws.Range("J7").Value = 424

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 47
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 94
$ws.Range("J3").Value = 92
$ws.Range("J7").Value = 287

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J6").Value = 89
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("J6").Value = 33
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 66

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J6").Value = 33
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("H4").Value = 11
$ws.Range("H7").Value = 211

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 124

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 83
$ws.Range("I4").Value = 14
$ws.Range("I7").Value = 338
$ws.Range("J7").Value = 225

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J6").Value = 96
$ws.Range("J7").Value = 255

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 98
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("J2").Value = 25
$ws.Range("J7").Value = 54

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J3").Value = 21
$ws.Range("J7").Value = 87

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 22
